# Apply the requested edits to the workbook.
$wb = $excel.ActiveWorkbook

# 1. Rename the "Include from ..." sheets to "Include #0" / "Include #1"
$wb.Worksheets.Item("Include from BRCBHPMTUSS").Name = "Include #0"
$wb.Worksheets.Item("Include from BRTabelaSUS").Name = "Include #1"

# 2. Update the Metadata sheet values
$ws = $wb.Worksheets.Item("Metadata")

# Version row: "0.1.0" -> "1.0.0"
$ws.Range("B3").Value = "1.0.0"

# Contact row: "No display for ContactDetail" -> "null (http://www.saude.gov.br)"
$ws.Range("B10").Value = "null (http://www.saude.gov.br)"
